$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the short-url column (B2:B7) - shared across all data rows
$ws.Range("B2:B7").Value = "jJG0k2"

# Row 7: country of origin changes from Cuba (CUB, id 46) to Uzbekistan (UZB, id 203)
# F7, N7, O7 hold digit-only text (shared-string) values; a plain .Value
# assignment would be auto-typed as a number by Excel, changing t="s" to a
# numeric cell and pulling in a new number-format style. Writing a TEXT()
# formula and then pasting-as-values keeps the cell text-typed with its
# original style untouched.
$ws.Range("F7").Formula = '=TEXT(203,"0")'
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false

$ws.Range("G7").Value = "Uzbekistan"
$ws.Range("H7").Value = "UZB"
$ws.Range("I7").Value = "UZB"

# Row 7: swap refugees / asylum_seekers counts
$ws.Range("N7").Formula = '=TEXT(0,"0")'
$ws.Range("N7").Copy() | Out-Null
$ws.Range("N7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false

$ws.Range("O7").Formula = '=TEXT(5,"0")'
$ws.Range("O7").Copy() | Out-Null
$ws.Range("O7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false
